$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "吉视传媒"
$ws.Range("C14").Value = "大元泵业"
$ws.Range("C15").Value = "合力泰"
$ws.Range("C16").Value = "大位科技"
$ws.Range("C17").Value = "岩山科技"
$ws.Range("C18").Value = "新易盛"
$ws.Range("C19").Value = "天融信"
$ws.Range("C21").Value = "川润股份"
